# Sales Order Module Creating Feature
# Adds a new "CreateSO" worksheet after the existing sheets, with a header
# row and sample login/sales-order test data.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet as the last tab -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CreateSO"

# --- 2. Pre-format the handful of cells whose literal text would otherwise
#        be auto-coerced to a number (leading zeros / plain digit strings).
#        Applying the Text format *before* the value is assigned preserves
#        the text exactly as typed.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"

# --- 3. Fill in all the cell values, strictly in row-major / left-to-right
#        order, so new shared strings are appended in the same order as the
#        original authoring session.

# Row 1 - header labels
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "ErrorMsg"
$ws.Range("E1").Value = "LoginPageTitle"
$ws.Range("F1").Value = "HomePageTitle"
$ws.Range("G1").Value = "SalesOrderPage"
$ws.Range("H1").Value = "Subject"
$ws.Range("I1").Value = "CustomerNumber"
$ws.Range("J1").Value = "CarrierName"
$ws.Range("K1").Value = "Salescommision"
$ws.Range("L1").Value = "Subject"
$ws.Range("M1").Value = "Due Date"

# Row 2 - sample data
$ws.Range("A2").Value = "ad123"
$ws.Range("B2").Value = 123456
$ws.Range("C2").Value = "PASS"
$ws.Range("D2").Value = "No such account configured for the user"
$ws.Range("E2").Value = "Zoho CRM - Sign in"
$ws.Range("F2").Value = "Zoho CRM - Home Page"
$ws.Range("G2").Value = "Zoho CRM - Displaying Custom View Details"
$ws.Range("H2").Value = "Hp laptop"
$ws.Range("I2").Value = "007"
$ws.Range("J2").Value = "BlueDart"
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = "Hp laptop"
$ws.Range("M2").Value = "'11/27/2019"

# Row 3
$ws.Range("A3").Value = "ad min"
$ws.Range("B3").Value = 123456
$ws.Range("C3").Value = "PASS"

# Row 4
$ws.Range("A4").Value = "ad$^%&"
$ws.Range("B4").Value = 123456
$ws.Range("C4").Value = "PASS"

# Row 5
$ws.Range("A5").Value = "ad123%^&"
$ws.Range("B5").Value = "123456"
$ws.Range("C5").Value = "PASS"

# Row 6
$ws.Range("A6").Value = 12345
$ws.Range("B6").Value = 123456
$ws.Range("C6").Value = "PASS"

# Row 7
$ws.Range("A7").Value = 45678
$ws.Range("B7").Value = 123456
$ws.Range("C7").Value = "PASS"

# --- 4. Apply the "Text" number format to the whole data block (rows 2-7)
#        This happens *after* the numeric values above are set, so numbers
#        like 123456 stay numeric even though the cell format is Text.
$ws.Range("A2:M2").NumberFormat = "@"
$ws.Range("A3:L3").NumberFormat = "@"
$ws.Range("A4:L4").NumberFormat = "@"
$ws.Range("A5:L5").NumberFormat = "@"
$ws.Range("A6:L6").NumberFormat = "@"
$ws.Range("A7:L7").NumberFormat = "@"

# --- 5. Header row styling - red font on yellow fill, Text number format
#        (column L keeps General format so it gets its own style entry).
$ws.Range("A1:K1,M1").NumberFormat = "@"
$ws.Range("A1:M1").Font.Color = 255
$ws.Range("A1:M1").Font.Size = 11
$ws.Range("A1:M1").Interior.Color = 65535

# --- 6. Sheet view / selection housekeeping to mirror the target workbook.
$ws.Range("A1:M7").Select()
$excel.ActiveWindow.Tabs.Item($ws.Name).Activate()
